# "Generate Report for Handback"
#
# The localization-status report is regenerated after a handback event for
# the 34f79766-1623-4b4a-8d87-f7b7d1f1ff8f.md file:
#   - Overview sheet: that file's zh-cn / de-de status flips from
#     "Ready for handoff" to "Handed back: in sync with en-US".
#   - zh-cn / de-de detail sheets: the row for that file gets its
#     "Latest Target File", "Latest Handback File" and
#     "Latest Handback DateTime" columns populated, and the new target-file
#     cell becomes a hyperlink (matching the existing hyperlink look used
#     elsewhere in the sheet).
#   - A few columns are widened so the newly-populated / longer text fits.

$wb = $excel.ActiveWorkbook

$srcMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/47cfc5a0aa3db05f2132aab8339b56635ca86443/e2e/34f79766-1623-4b4a-8d87-f7b7d1f1ff8f.md"

# ---------------------------------------------------------------------
# Overview sheet: row 2 is the 34f79766-....md file. It has been handed
# back and is now in sync with en-US (columns E = zh-cn, F = de-de).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# Row 3 (dca1bbfd-....md) keeps reporting "Ready for handoff" - rewrite it
# explicitly too so the shared-string table matches a freshly generated
# report rather than relying on incidental reuse.
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-09-04 10:43:50"

$overview.Columns.Item(5).ColumnWidth = 29.16666666666667
$overview.Columns.Item(6).ColumnWidth = 29.16666666666667

# ---------------------------------------------------------------------
# Helper: populate the handback columns for the 34f79766-....md row (row 2)
# on a locale detail sheet, and widen the columns that now hold longer text.
# ---------------------------------------------------------------------
function Update-LocaleSheet($sheetName, $xliffName, $handbackDatetime) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("C2").Value = "Ready for handoff"

    $target = $ws.Range("I2")
    $target.Value = "34f79766-1623-4b4a-8d87-f7b7d1f1ff8f.md"
    $target.Font.Underline = $true
    $target.Font.Color = 15570276
    $ws.Hyperlinks.Add($target, $srcMdUrl, [Type]::Missing, [Type]::Missing, "34f79766-1623-4b4a-8d87-f7b7d1f1ff8f.md")

    $ws.Range("J2").Value = $xliffName
    $ws.Range("K2").Value = $handbackDatetime

    $ws.Columns.Item(3).ColumnWidth = 29.16666666666667
    $ws.Columns.Item(9).ColumnWidth = 39.16666666666667
    $ws.Columns.Item(10).ColumnWidth = 39.16666666666667
}

Update-LocaleSheet "zh-cn" "34f79766-1623-4b4a-8d87-f7b7d1f1ff8f.eaec23e3727f30c0a17a3321392160a545743c7c.zh-cn.xlf" "2016-09-04 10:44:30"
Update-LocaleSheet "de-de" "34f79766-1623-4b4a-8d87-f7b7d1f1ff8f.eaec23e3727f30c0a17a3321392160a545743c7c.de-de.xlf" "2016-09-04 10:44:37"

Write-Output "Generate Report for Handback: applied"
